$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text so numeric-looking strings (e.g. "227.75") are not
# auto-converted to numbers by Excel -- the source workbook stores the whole
# Price column as literal text (note some values use "." as a thousands
# separator, e.g. "36.047.86", which would not even parse as a number).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '36.047.86'
$ws.Cells.Item(2, 5).Value = '  -3.80%  '
$ws.Cells.Item(3, 4).Value = '1.945.41'
$ws.Cells.Item(3, 5).Value = '  -3.79%  '
$ws.Cells.Item(4, 5).Value = '  +0.11%  '
$ws.Cells.Item(5, 4).Value = '227.75'
$ws.Cells.Item(5, 5).Value = '  -10.16%  '
$ws.Cells.Item(6, 5).Value = '  -4.75%  '
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 4).Value = '52.51'
$ws.Cells.Item(8, 5).Value = '  -7.16%  '
$ws.Cells.Item(9, 4).Value = '0.362'
$ws.Cells.Item(9, 5).Value = '  -4.95%  '
$ws.Cells.Item(10, 4).Value = '56.62'
$ws.Cells.Item(10, 5).Value = '  -0.80%  '
$ws.Cells.Item(11, 4).Value = '0.0726'
$ws.Cells.Item(11, 5).Value = '  -6.86%  '
$ws.Cells.Item(12, 5).Value = '  -4.58%  '
$ws.Cells.Item(13, 4).Value = '2.235.18'
$ws.Cells.Item(13, 5).Value = '  -3.68%  '
$ws.Cells.Item(14, 4).Value = '13.66'
$ws.Cells.Item(14, 5).Value = '  -5.61%  '
$ws.Cells.Item(15, 4).Value = '0.736'
$ws.Cells.Item(15, 5).Value = '  -9.20%  '
$ws.Cells.Item(16, 4).Value = '19.22'
$ws.Cells.Item(16, 5).Value = '  -8.08%  '
$ws.Cells.Item(17, 4).Value = '1.958.23'
$ws.Cells.Item(17, 5).Value = '  -3.62%  '
$ws.Cells.Item(18, 4).Value = '4.92'
$ws.Cells.Item(18, 5).Value = '  -7.34%  '
$ws.Cells.Item(19, 4).Value = '35.986.84'
$ws.Cells.Item(19, 5).Value = '  -3.72%  '
$ws.Cells.Item(20, 4).Value = '66.65'
$ws.Cells.Item(20, 5).Value = '  -3.99%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0781'
$ws.Cells.Item(21, 5).Value = '  -7.59%  '
$ws.Cells.Item(22, 4).Value = '4.89'
$ws.Cells.Item(22, 5).Value = '  -5.40%  '
$ws.Cells.Item(23, 4).Value = '219.41'
$ws.Cells.Item(23, 5).Value = '  -3.72%  '
$ws.Cells.Item(24, 5).Value = '  -0.04%  '
$ws.Cells.Item(25, 4).Value = '2.32'
$ws.Cells.Item(25, 5).Value = '  -0.09%  '
$ws.Cells.Item(26, 5).Value = '  -11.85%  '
$ws.Cells.Item(27, 4).Value = '159.84'
$ws.Cells.Item(27, 5).Value = '  -2.04%  '
$ws.Cells.Item(28, 4).Value = '8.31'
$ws.Cells.Item(28, 5).Value = '  -7.45%  '
$ws.Cells.Item(29, 4).Value = '18.52'
$ws.Cells.Item(29, 5).Value = '  -6.10%  '
$ws.Cells.Item(30, 4).Value = '1.28'
$ws.Cells.Item(30, 5).Value = '  -7.04%  '
$ws.Cells.Item(31, 5).Value = '  -11.27%  '
$ws.Cells.Item(32, 5).Value = '  -4.38%  '
$ws.Cells.Item(33, 5).Value = '  -9.16%  '
$ws.Cells.Item(34, 5).Value = '  -10.98%  '
$ws.Cells.Item(35, 4).Value = '4.13'
$ws.Cells.Item(35, 5).Value = '  -8.97%  '
$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).Value = '2.23'
$ws.Cells.Item(36, 5).Value = '  -7.78%  '
$ws.Cells.Item(37, 2).Value = 'BinanceUSD'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(37, 4).Value = '1.00'
$ws.Cells.Item(37, 5).Value = '  -0.17%  '
$ws.Cells.Item(38, 5).Value = '  -2.45%  '
$ws.Cells.Item(39, 4).Value = '3.09'
$ws.Cells.Item(39, 5).Value = '  -8.81%  '
$ws.Cells.Item(40, 5).Value = '  -1.52%  '
$ws.Cells.Item(41, 4).Value = '4.94'
$ws.Cells.Item(41, 5).Value = '  -6.77%  '
$ws.Cells.Item(42, 4).Value = '1.387.89'
$ws.Cells.Item(42, 5).Value = '  -1.73%  '
$ws.Cells.Item(43, 5).Value = '  -8.77%  '
$ws.Cells.Item(44, 4).Value = '0.0846'
$ws.Cells.Item(44, 5).Value = '  -11.91%  '
$ws.Cells.Item(45, 4).Value = '1.05'
$ws.Cells.Item(45, 5).Value = '  -12.78%  '
$ws.Cells.Item(46, 4).Value = '85.33'
$ws.Cells.Item(46, 5).Value = '  -5.49%  '
$ws.Cells.Item(47, 4).Value = '0.962'
$ws.Cells.Item(47, 5).Value = '  -6.50%  '
$ws.Cells.Item(48, 2).Value = 'MXToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(48, 4).Value = '2.83'
$ws.Cells.Item(48, 5).Value = '  -1.27%  '
$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).Value = '14.38'
$ws.Cells.Item(49, 5).Value = '  -9.44%  '
$ws.Cells.Item(50, 4).Value = '6.64'
$ws.Cells.Item(50, 5).Value = '  -8.78%  '
$ws.Cells.Item(51, 4).Value = '2.129.24'

# Restore the default style on column D (clears the quirky "s" style index
# that NumberFormat assignment introduces) so only cell VALUES differ from
# the original workbook, matching the source diff.
$ws.Range("D2:D51").Style = "Normal"
